$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3224.125
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 3224.125
$ws.Range("N88").Value = -4036.125

$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3224.125
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 3224.125
$ws.Range("N91").Value = -6032.125

$ws.Range("H98").Value = 4460.6924
$ws.Range("I98").Value = 4338.6
$ws.Range("J98").Value = 4537
$ws.Range("K98").Value = 4338.6
$ws.Range("L98").Value = 4537
$ws.Range("M98").Value = -2840.6
$ws.Range("N98").Value = -7533

$ws.Range("H106").Value = 8459.4
$ws.Range("I106").Value = 8459.4
$ws.Range("K106").Value = 8459.4
$ws.Range("M106").Value = -7828.4

$ws.Range("H116").Value = 5278
$ws.Range("I116").Value = 5278
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 5278
$ws.Range("L116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -1836

$ws.Range("H122").Value = 4460.6924
$ws.Range("I122").Value = 4338.6
$ws.Range("J122").Value = 4537
$ws.Range("K122").Value = 13015.8
$ws.Range("L122").Value = 13611
$ws.Range("M122").Value = -10565.8
$ws.Range("N122").Value = -18511

$ws.Range("H125").Value = 1749.75
$ws.Range("I125").Value = 1499.5
$ws.Range("K125").Value = 13495.5
$ws.Range("M125").Value = -11035.5

$ws.Range("H132").Value = 2882.2856
$ws.Range("I132").Value = 2908.2354
$ws.Range("K132").Value = 8724.706200000001
$ws.Range("M132").Value = -6194.706200000001

$ws.Range("H138").Value = 2604.492
$ws.Range("J138").Value = 3398.5676
$ws.Range("L138").Value = 10195.7028
$ws.Range("N138").Value = -20475.7028

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2765.3333
$ws.Range("I32").Value = 1163.6167
$ws.Range("K32").Value = 1163.6167
$ws.Range("M32").Value = -876.6167

$ws.Range("H74").Value = 1436.2727
$ws.Range("I74").Value = 1496.8
$ws.Range("K74").Value = 1496.8
$ws.Range("M74").Value = -622.8

$ws.Range("H77").Value = 1436.2727
$ws.Range("I77").Value = 1496.8
$ws.Range("K77").Value = 7484
$ws.Range("M77").Value = -3116

$ws.Range("H88").Value = 3721.25
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 3721.25
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").Value = 3721.25
$ws.Range("N88").Value = -4533.25

$ws.Range("H91").Value = 3721.25
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 3721.25
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").Value = 3721.25
$ws.Range("N91").Value = -6529.25

$ws.Range("H122").Value = 3660.1428
$ws.Range("I122").Value = 3337.1667
$ws.Range("K122").Value = 10011.5001
$ws.Range("M122").Value = -7561.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3748.5334
$ws.Range("I20").Value = 3558.1
$ws.Range("J20").Value = 4129.4
$ws.Range("K20").Value = 3558.1
$ws.Range("L20").Value = 4129.4
$ws.Range("M20").Value = -3311.1
$ws.Range("N20").Value = -4623.4

$ws.Range("H63").Value = 39635
$ws.Range("J63").Value = 39635
$ws.Range("L63").Value = 39635
$ws.Range("N63").Value = -41007

$ws.Range("H66").Value = 39635
$ws.Range("J66").Value = 39635
$ws.Range("L66").Value = 118905
$ws.Range("N66").Value = -125769

$ws.Range("H86").Value = 23615684
$ws.Range("I86").Value = 70835336
$ws.Range("J86").Value = 5859.25
$ws.Range("K86").Value = 70835336
$ws.Range("L86").Value = 5859.25
$ws.Range("M86").Value = -70834213
$ws.Range("N86").Value = -8105.25

$ws.Range("H89").Value = 23615684
$ws.Range("I89").Value = 70835336
$ws.Range("J89").Value = 5859.25
$ws.Range("K89").Value = 354176680
$ws.Range("L89").Value = 29296.25
$ws.Range("M89").Value = -354171064
$ws.Range("N89").Value = -40528.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 7164028.5

$ws.Range("H86").Value = 16908.545
$ws.Range("I86").Value = 23082.334
$ws.Range("K86").Value = 23082.334
$ws.Range("M86").Value = -21959.334

$ws.Range("H89").Value = 16908.545
$ws.Range("I89").Value = 23082.334
$ws.Range("K89").Value = 115411.67
$ws.Range("M89").Value = -109795.67

$ws.Range("H99").Value = 10611.429
$ws.Range("I99").Value = 6595.3335
$ws.Range("K99").Value = 6595.3335
$ws.Range("M99").Value = -5097.3335

$ws.Range("H126").Value = 10611.429
$ws.Range("I126").Value = 6595.3335
$ws.Range("K126").Value = 19786.0005
$ws.Range("M126").Value = -17316.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 4042.25
$ws.Range("J2").Value = 570
$ws.Range("L2").Value = 3420
$ws.Range("N2").Value = -3646

$ws.Range("H11").Value = 1441.2142
$ws.Range("I11").Value = 575.3333
$ws.Range("K11").Value = 1725.9999
$ws.Range("M11").Value = -1585.9999

$ws.Range("H26").Value = 1226.8
$ws.Range("I26").Value = 1733.5714
$ws.Range("J26").Value = 44.333332
$ws.Range("K26").Value = 5200.7142
$ws.Range("L26").Value = 132.999996
$ws.Range("M26").Value = -4912.7142
$ws.Range("N26").Value = -708.999996

$ws.Range("H37").Value = 69960.336
$ws.Range("J37").Value = 69960.336
$ws.Range("L37").Value = 209881.008
$ws.Range("N37").Value = -210105.008

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 10001.667
$ws.Range("J5").Value = 10005
$ws.Range("L5").Value = 10005
$ws.Range("N5").Value = -10229

$ws.Range("H7").Value = 85841680
$ws.Range("I7").Value = 200040.6
$ws.Range("J7").Value = 147014290
$ws.Range("K7").Value = 200040.6
$ws.Range("L7").Value = 147014290
$ws.Range("M7").Value = -199928.6
$ws.Range("N7").Value = -147014514

$ws.Range("H8").Value = 85841680
$ws.Range("I8").Value = 200040.6
$ws.Range("J8").Value = 147014290
$ws.Range("K8").Value = 200040.6
$ws.Range("L8").Value = 147014290
$ws.Range("M8").Value = -199901.6
$ws.Range("N8").Value = -147014568

$ws.Range("H62").Value = 48900
$ws.Range("J62").Value = 48900
$ws.Range("L62").Value = 48900
$ws.Range("N62").Value = -50272

$ws.Range("H65").Value = 48900
$ws.Range("J65").Value = 48900
$ws.Range("L65").Value = 146700
$ws.Range("N65").Value = -153564

$ws.Range("H80").Value = 7461.2856
$ws.Range("I80").Value = 6848.4
$ws.Range("K80").Value = 6848.4
$ws.Range("M80").Value = -5850.4

$ws.Range("H83").Value = 7461.2856
$ws.Range("I83").Value = 6848.4
$ws.Range("K83").Value = 34242
$ws.Range("M83").Value = -29250

$ws.Range("H132").Value = 6760.205
$ws.Range("I132").Value = 6790.5
$ws.Range("K132").Value = 20371.5
$ws.Range("M132").Value = -17841.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 26871.572
$ws.Range("I7").Value = 26871.572
$ws.Range("K7").Value = 26871.572
$ws.Range("M7").Value = -26759.572

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").ClearContents()
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = 0

$ws.Range("H93").Value = 1604.625
$ws.Range("I93").Value = 1639.6666
$ws.Range("J93").Value = 1499.5
$ws.Range("K93").Value = 1639.6666
$ws.Range("L93").Value = 1499.5
$ws.Range("M93").Value = -391.6666
$ws.Range("N93").Value = -3995.5

$ws.Range("H100").Value = 2016.3334
$ws.Range("I100").Value = 1873.4375
$ws.Range("J100").Value = 2302.125
$ws.Range("K100").Value = 1873.4375
$ws.Range("L100").Value = 2302.125
$ws.Range("M100").Value = -1332.4375
$ws.Range("N100").Value = -3384.125

$ws.Range("H126").Value = 26871.572
$ws.Range("I126").Value = 26871.572
$ws.Range("K126").Value = 80614.716
$ws.Range("M126").Value = -78144.716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 16696.5
$ws.Range("I62").Value = 16696.5
$ws.Range("K62").Value = 16696.5
$ws.Range("M62").Value = -16072.5

$ws.Range("H65").Value = 16696.5
$ws.Range("I65").Value = 16696.5
$ws.Range("K65").Value = 83482.5
$ws.Range("M65").Value = -80362.5

$ws.Range("H81").Value = 7099.3335
$ws.Range("J81").Value = 1298
$ws.Range("L81").Value = 2596
$ws.Range("N81").Value = -4718

$ws.Range("H84").Value = 7099.3335
$ws.Range("J84").Value = 1298
$ws.Range("L84").Value = 12980
$ws.Range("N84").Value = -23588

$ws.Range("H126").Value = 1998
$ws.Range("I126").Value = 1998
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5994
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -3524
